# Auto-applies numeric corrections to the Ragnarok_Profits workbook sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) per the commit's
# scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 9910.200000000001
$ws.Range("J69").Value = 9910.200000000001
$ws.Range("L69").Value = 29730.6
$ws.Range("N69").Value = -31478.6
$ws.Range("H72").Value = 9910.200000000001
$ws.Range("J72").Value = 9910.200000000001
$ws.Range("L72").Value = 89191.8
$ws.Range("N72").Value = -97927.8
$ws.Range("H95").Value = 79998
$ws.Range("J95").Value = 79998
$ws.Range("L95").Value = 79998
$ws.Range("N95").Value = -85490
$ws.Range("H100").Value = 7208.5835
$ws.Range("I100").Value = 6103.3335
$ws.Range("J100").Value = 8313.833000000001
$ws.Range("K100").Value = 6103.3335
$ws.Range("L100").Value = 8313.833000000001
$ws.Range("M100").Value = -5562.3335
$ws.Range("N100").Value = -9395.833000000001
$ws.Range("H115").Value = 593.625
$ws.Range("I115").Value = 393.57144
$ws.Range("K115").Value = 1180.71432
$ws.Range("M115").Value = 386.28568
$ws.Range("H132").Value = 1718.0834
$ws.Range("I132").Value = 1665.9556
$ws.Range("K132").Value = 4997.8668
$ws.Range("M132").Value = -2467.8668
$ws.Range("H135").Value = 2444.625
$ws.Range("I135").Value = 311.6
$ws.Range("K135").Value = 2804.4
$ws.Range("M135").Value = -269.4000000000001
$ws.Range("H138").Value = 4001.3242
$ws.Range("I138").Value = 2749.1667
$ws.Range("J138").Value = 5187.579
$ws.Range("K138").Value = 8247.500100000001
$ws.Range("L138").Value = 15562.737
$ws.Range("M138").Value = -3107.500100000001
$ws.Range("N138").Value = -25842.737

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 5251.8423
$ws.Range("I110").Value = 4585.933
$ws.Range("J110").Value = 7749
$ws.Range("K110").Value = 4585.933
$ws.Range("L110").Value = 7749
$ws.Range("M110").Value = -2540.933
$ws.Range("N110").Value = -11839
$ws.Range("H132").Value = 2858774.2
$ws.Range("I132").Value = 1573.4
$ws.Range("K132").Value = 4720.200000000001
$ws.Range("M132").Value = -2190.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 37655.75
$ws.Range("J95").Value = 37655.75
$ws.Range("L95").Value = 37655.75
$ws.Range("N95").Value = -43147.75
$ws.Range("H100").Value = 56398
$ws.Range("J100").Value = 56398
$ws.Range("L100").Value = 56398
$ws.Range("N100").Value = -58562
$ws.Range("H132").Value = 109500
$ws.Range("J132").Value = 109500
$ws.Range("L132").Value = 109500
$ws.Range("N132").Value = -119620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 83666.664
$ws.Range("J28").Value = 83666.664
$ws.Range("L28").Value = 83666.664
$ws.Range("N28").Value = -84156.664
$ws.Range("H107").Value = 1167.7826
$ws.Range("I107").Value = 288.85715
$ws.Range("J107").Value = 2535
$ws.Range("K107").Value = 288.85715
$ws.Range("L107").Value = 2535
$ws.Range("M107").Value = 1631.14285
$ws.Range("N107").Value = -6375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 45457240
$ws.Range("J12").Value = 2642.7368
$ws.Range("L12").Value = 7928.2104
$ws.Range("N12").Value = -8274.2104
$ws.Range("H17").Value = 7752
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 9242.4
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 27727.2
$ws.Range("M17").Value = -731
$ws.Range("N17").Value = -28065.2
$ws.Range("H63").Value = 26384.4
$ws.Range("I63").Value = 10256
$ws.Range("K63").Value = 30768
$ws.Range("M63").Value = -30019
$ws.Range("H66").Value = 26384.4
$ws.Range("I66").Value = 10256
$ws.Range("K66").Value = 92304
$ws.Range("M66").Value = -88560
$ws.Range("H75").Value = 33332
$ws.Range("J75").Value = 33332
$ws.Range("L75").Value = 99996
$ws.Range("N75").Value = -101992
$ws.Range("H78").Value = 33332
$ws.Range("J78").Value = 33332
$ws.Range("L78").Value = 299988
$ws.Range("N78").Value = -309972
$ws.Range("H114").Value = 33333
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 33333
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 99999
$ws.Range("M114").ClearContents()
$ws.Range("N114").Value = -106507

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2857.8572
$ws.Range("J80").Value = 3201.5
$ws.Range("L80").Value = 3201.5
$ws.Range("N80").Value = -5197.5
$ws.Range("H83").Value = 2857.8572
$ws.Range("J83").Value = 3201.5
$ws.Range("L83").Value = 16007.5
$ws.Range("N83").Value = -25991.5
$ws.Range("H102").Value = 1775.5834
$ws.Range("I102").Value = 1787.5652
$ws.Range("K102").Value = 1787.5652
$ws.Range("M102").Value = -165.5652
$ws.Range("H126").Value = 6068.222
$ws.Range("I126").Value = 5459.3335
$ws.Range("J126").Value = 7286
$ws.Range("K126").Value = 16378.0005
$ws.Range("L126").Value = 21858
$ws.Range("M126").Value = -13908.0005
$ws.Range("N126").Value = -26798

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8749.357
$ws.Range("I7").Value = 8541.083000000001
$ws.Range("J7").Value = 9999
$ws.Range("K7").Value = 8541.083000000001
$ws.Range("L7").Value = 9999
$ws.Range("M7").Value = -8429.083000000001
$ws.Range("N7").Value = -10223
$ws.Range("H22").Value = 10585.333
$ws.Range("J22").Value = 3636.8462
$ws.Range("L22").Value = 3636.8462
$ws.Range("N22").Value = -4226.8462
$ws.Range("H27").Value = 10585.333
$ws.Range("J27").Value = 3636.8462
$ws.Range("L27").Value = 3636.8462
$ws.Range("N27").Value = -3850.8462
$ws.Range("H82").Value = 4443.1113
$ws.Range("I82").Value = 1165
$ws.Range("K82").Value = 1165
$ws.Range("M82").Value = -804
$ws.Range("H85").Value = 4443.1113
$ws.Range("I85").Value = 1165
$ws.Range("K85").Value = 1165
$ws.Range("M85").Value = 83
$ws.Range("H100").Value = 25029410
$ws.Range("I100").Value = 2930.4
$ws.Range("K100").Value = 2930.4
$ws.Range("M100").Value = -2389.4
$ws.Range("H126").Value = 8749.357
$ws.Range("I126").Value = 8541.083000000001
$ws.Range("J126").Value = 9999
$ws.Range("K126").Value = 25623.249
$ws.Range("L126").Value = 29997
$ws.Range("M126").Value = -23153.249
$ws.Range("N126").Value = -34937
$ws.Range("H129").Value = 90429
$ws.Range("J129").Value = 90429
$ws.Range("L129").Value = 90429
$ws.Range("N129").Value = -100429
$ws.Range("H132").Value = 3913.578
$ws.Range("I132").Value = 2525.625
$ws.Range("J132").Value = 5499.8096
$ws.Range("K132").Value = 7576.875
$ws.Range("L132").Value = 16499.4288
$ws.Range("M132").Value = -5046.875
$ws.Range("N132").Value = -21559.4288

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 22246
$ws.Range("I45").Value = 48000
$ws.Range("J45").Value = 13661.333
$ws.Range("K45").Value = 48000
$ws.Range("L45").Value = 13661.333
$ws.Range("M45").Value = -47509
$ws.Range("N45").Value = -14643.333
$ws.Range("H62").Value = 9647.375
$ws.Range("I62").Value = 4070
$ws.Range("J62").Value = 15224.75
$ws.Range("K62").Value = 4070
$ws.Range("L62").Value = 15224.75
$ws.Range("M62").Value = -3446
$ws.Range("N62").Value = -16472.75
$ws.Range("H65").Value = 9647.375
$ws.Range("I65").Value = 4070
$ws.Range("J65").Value = 15224.75
$ws.Range("K65").Value = 20350
$ws.Range("L65").Value = 76123.75
$ws.Range("M65").Value = -17230
$ws.Range("N65").Value = -82363.75
$ws.Range("H96").Value = 11344.75
$ws.Range("I96").Value = 11153.8
$ws.Range("J96").Value = 11663
$ws.Range("K96").Value = 11153.8
$ws.Range("L96").Value = 11663
$ws.Range("M96").Value = -9780.799999999999
$ws.Range("N96").Value = -14409
$ws.Range("H122").Value = 2483.25
$ws.Range("I122").Value = 2114.4285
$ws.Range("K122").Value = 6343.2855
$ws.Range("M122").Value = -3893.2855
$ws.Range("H132").Value = 224602.25
$ws.Range("I132").Value = 2495.2368
$ws.Range("K132").Value = 7485.7104
$ws.Range("M132").Value = -4955.7104

